# Add two new puzzle-solver result rows to the "Five Puzzles" sheet,
# matching the formatting of the existing data row above them.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Five Puzzles")

# Copy formatting (cell styles) from the last existing data row (row 4)
# down into the two new rows (5 and 6) before filling in values.
$ws.Range("A4:DI4").Copy()
$ws.Range("A5:DI5").PasteSpecial(-4122)
$ws.Range("A4:DI4").Copy()
$ws.Range("A6:DI6").PasteSpecial(-4122)

# Fill in the new data for rows 5 and 6.
    $ws.Range("A5").Value = "multipuzzle"
    $ws.Range("B5").Value = 2
    $ws.Range("C5").Value = "pomeranz_805_1 | pomeranz_805_2 | pomeranz_805_8 | pomeranz_805_9 | pomeranz_805_10"
    $ws.Range("D5").Value = 5
    $ws.Range("E5").Value = 5
    $ws.Range("F5").Value = 0
    $ws.Range("G5").Value = "Input Puzzle #0"
    $ws.Range("H5").Value = 805
    $ws.Range("I5").Value = "Solved Puzzle #0"
    $ws.Range("J5").Value = 811
    $ws.Range("K5").Value = 0
    $ws.Range("L5").Value = 0.0073982737361300002
    $ws.Range("M5").Value = 0.93958076448799999
    $ws.Range("N5").Value = 0
    $ws.Range("O5").Value = 43
    $ws.Range("P5").Value = "Solved Puzzle #0"
    $ws.Range("Q5").Value = 811
    $ws.Range("R5").Value = 0.72996300863100005
    $ws.Range("S5").Value = 0.0073982737361300002
    $ws.Range("T5").Value = 0.20961775585699999
    $ws.Range("U5").Value = 0
    $ws.Range("V5").Value = 43
    $ws.Range("W5").Value = "Solved Puzzle#0"
    $ws.Range("X5").Value = 786
    $ws.Range("Y5").Value = 0.75030156815399995
    $ws.Range("Z5").Value = 43
    $ws.Range("AA5").Value = 24
    $ws.Range("AB5").Value = "Input Puzzle #1"
    $ws.Range("AC5").Value = 805
    $ws.Range("AD5").Value = "Solved Puzzle #2"
    $ws.Range("AE5").Value = 841
    $ws.Range("AF5").Value = 0
    $ws.Range("AG5").Value = 0.042806183115300003
    $ws.Range("AH5").Value = 0.92508917954799996
    $ws.Range("AI5").Value = 0.0011890606420899999
    $ws.Range("AJ5").Value = 26
    $ws.Range("AK5").Value = "Solved Puzzle #2"
    $ws.Range("AL5").Value = 841
    $ws.Range("AM5").Value = 0
    $ws.Range("AN5").Value = 0.042806183115300003
    $ws.Range("AO5").Value = 0.92508917954799996
    $ws.Range("AP5").Value = 0.0011890606420899999
    $ws.Range("AQ5").Value = 26
    $ws.Range("AR5").Value = "Solved Puzzle#2"
    $ws.Range("AS5").Value = 923
    $ws.Range("AT5").Value = 0.72655426764999997
    $ws.Range("AU5").Value = 26
    $ws.Range("AV5").Value = 144
    $ws.Range("AW5").Value = "Input Puzzle #2"
    $ws.Range("AX5").Value = 805
    $ws.Range("AY5").Value = "Solved Puzzle #1"
    $ws.Range("AZ5").Value = 814
    $ws.Range("BA5").Value = 0
    $ws.Range("BB5").Value = 0.011056511056500001
    $ws.Range("BC5").Value = 0.98894348894299999
    $ws.Range("BD5").Value = 0
    $ws.Range("BE5").Value = 0
    $ws.Range("BF5").Value = "Solved Puzzle #1"
    $ws.Range("BG5").Value = 814
    $ws.Range("BH5").Value = 0.98894348894299999
    $ws.Range("BI5").Value = 0.011056511056500001
    $ws.Range("BJ5").Value = 0
    $ws.Range("BK5").Value = 0
    $ws.Range("BL5").Value = 0
    $ws.Range("BM5").Value = "Solved Puzzle#1"
    $ws.Range("BN5").Value = 841
    $ws.Range("BO5").Value = 0.95451843043999995
    $ws.Range("BP5").Value = 0
    $ws.Range("BQ5").Value = 36
    $ws.Range("BR5").Value = "Input Puzzle #3"
    $ws.Range("BS5").Value = 805
    $ws.Range("BT5").Value = "Solved Puzzle #4"
    $ws.Range("BU5").Value = 807
    $ws.Range("BV5").Value = 0
    $ws.Range("BW5").Value = 0.0024783147459699999
    $ws.Range("BX5").Value = 0.99504337050799996
    $ws.Range("BY5").Value = 0
    $ws.Range("BZ5").Value = 2
    $ws.Range("CA5").Value = "Solved Puzzle #4"
    $ws.Range("CB5").Value = 807
    $ws.Range("CC5").Value = 0.79553903345700006
    $ws.Range("CD5").Value = 0.0024783147459699999
    $ws.Range("CE5").Value = 0.198265179678
    $ws.Range("CF5").Value = 0.0012391573729900001
    $ws.Range("CG5").Value = 2
    $ws.Range("CH5").Value = "Solved Puzzle#4"
    $ws.Range("CI5").Value = 811
    $ws.Range("CJ5").Value = 0.86961869618700005
    $ws.Range("CK5").Value = 2
    $ws.Range("CL5").Value = 8
    $ws.Range("CM5").Value = "Input Puzzle #4"
    $ws.Range("CN5").Value = 805
    $ws.Range("CO5").Value = "Solved Puzzle #3"
    $ws.Range("CP5").Value = 823
    $ws.Range("CQ5").Value = 0
    $ws.Range("CR5").Value = 0.021871202916200001
    $ws.Range("CS5").Value = 0.97812879708400002
    $ws.Range("CT5").Value = 0
    $ws.Range("CU5").Value = 0
    $ws.Range("CV5").Value = "Solved Puzzle #3"
    $ws.Range("CW5").Value = 823
    $ws.Range("CX5").Value = 0.97812879708400002
    $ws.Range("CY5").Value = 0.021871202916200001
    $ws.Range("CZ5").Value = 0
    $ws.Range("DA5").Value = 0
    $ws.Range("DB5").Value = 0
    $ws.Range("DC5").Value = "Solved Puzzle#3"
    $ws.Range("DD5").Value = 877
    $ws.Range("DE5").Value = 0.91277080957800005
    $ws.Range("DF5").Value = 0
    $ws.Range("DG5").Value = 72
    $ws.Range("A6").Value = "paikin_tal"
    $ws.Range("B6").Value = 2
    $ws.Range("C6").Value = "pomeranz_805_1 | pomeranz_805_2 | pomeranz_805_8 | pomeranz_805_9 | pomeranz_805_10"
    $ws.Range("D6").Value = 5
    $ws.Range("E6").Value = 5
    $ws.Range("F6").Value = 0
    $ws.Range("G6").Value = "Input Puzzle #0"
    $ws.Range("H6").Value = 805
    $ws.Range("I6").Value = "Solved Puzzle #1"
    $ws.Range("J6").Value = 808
    $ws.Range("K6").Value = 0
    $ws.Range("L6").Value = 0.0037128712871299999
    $ws.Range("M6").Value = 0.67450495049500003
    $ws.Range("N6").Value = 0
    $ws.Range("O6").Value = 260
    $ws.Range("P6").Value = "Solved Puzzle #1"
    $ws.Range("Q6").Value = 808
    $ws.Range("R6").Value = 0.48886138613899999
    $ws.Range("S6").Value = 0.0037128712871299999
    $ws.Range("T6").Value = 0.18564356435599999
    $ws.Range("U6").Value = 0
    $ws.Range("V6").Value = 260
    $ws.Range("W6").Value = "Solved Puzzle#1"
    $ws.Range("X6").Value = 557
    $ws.Range("Y6").Value = 0.50642594859199996
    $ws.Range("Z6").Value = 260
    $ws.Range("AA6").Value = 12
    $ws.Range("AB6").Value = "Input Puzzle #1"
    $ws.Range("AC6").Value = 805
    $ws.Range("AD6").Value = "Solved Puzzle #4"
    $ws.Range("AE6").Value = 1794
    $ws.Range("AF6").Value = 0
    $ws.Range("AG6").Value = 0.55128205128200003
    $ws.Range("AH6").Value = 0.43255295429200002
    $ws.Range("AI6").Value = 0
    $ws.Range("AJ6").Value = 29
    $ws.Range("AK6").Value = "Solved Puzzle #4"
    $ws.Range("AL6").Value = 1794
    $ws.Range("AM6").Value = 0
    $ws.Range("AN6").Value = 0.55128205128200003
    $ws.Range("AO6").Value = 0.43255295429200002
    $ws.Range("AP6").Value = 0
    $ws.Range("AQ6").Value = 29
    $ws.Range("AR6").Value = "Solved Puzzle#4"
    $ws.Range("AS6").Value = 4732
    $ws.Range("AT6").Value = 0.13820625918900001
    $ws.Range("AU6").Value = 29
    $ws.Range("AV6").Value = 3956
    $ws.Range("AW6").Value = "Input Puzzle #2"
    $ws.Range("AX6").Value = 805
    $ws.Range("AY6").Value = "Solved Puzzle #2"
    $ws.Range("AZ6").Value = 1609
    $ws.Range("BA6").Value = 0
    $ws.Range("BB6").Value = 0.49968924797999997
    $ws.Range("BC6").Value = 0.50031075201999997
    $ws.Range("BD6").Value = 0
    $ws.Range("BE6").Value = 0
    $ws.Range("BF6").Value = "Solved Puzzle #2"
    $ws.Range("BG6").Value = 1609
    $ws.Range("BH6").Value = 0.50031075201999997
    $ws.Range("BI6").Value = 0.49968924797999997
    $ws.Range("BJ6").Value = 0
    $ws.Range("BK6").Value = 0
    $ws.Range("BL6").Value = 0
    $ws.Range("BM6").Value = "Solved Puzzle#2"
    $ws.Range("BN6").Value = 4021
    $ws.Range("BO6").Value = 0.198706789356
    $ws.Range("BP6").Value = 0
    $ws.Range("BQ6").Value = 3216
    $ws.Range("BR6").Value = "Input Puzzle #3"
    $ws.Range("BS6").Value = 805
    $ws.Range("BT6").Value = "Solved Puzzle #2"
    $ws.Range("BU6").Value = 1611
    $ws.Range("BV6").Value = 0
    $ws.Range("BW6").Value = 0.50031036623199998
    $ws.Range("BX6").Value = 0.49844816883900001
    $ws.Range("BY6").Value = 0
    $ws.Range("BZ6").Value = 2
    $ws.Range("CA6").Value = "Solved Puzzle #2"
    $ws.Range("CB6").Value = 1611
    $ws.Range("CC6").Value = 0
    $ws.Range("CD6").Value = 0.50031036623199998
    $ws.Range("CE6").Value = 0.49844816883900001
    $ws.Range("CF6").Value = 0
    $ws.Range("CG6").Value = 2
    $ws.Range("CH6").Value = "Solved Puzzle#2"
    $ws.Range("CI6").Value = 4027
    $ws.Range("CJ6").Value = 0.17386448250200001
    $ws.Range("CK6").Value = 2
    $ws.Range("CL6").Value = 3224
    $ws.Range("CM6").Value = "Input Puzzle #4"
    $ws.Range("CN6").Value = 805
    $ws.Range("CO6").Value = "Solved Puzzle #4"
    $ws.Range("CP6").Value = 1766
    $ws.Range("CQ6").Value = 0
    $ws.Range("CR6").Value = 0.54416761041899997
    $ws.Range("CS6").Value = 0.455266138165
    $ws.Range("CT6").Value = 0
    $ws.Range("CU6").Value = 1
    $ws.Range("CV6").Value = "Solved Puzzle #4"
    $ws.Range("CW6").Value = 1766
    $ws.Range("CX6").Value = 0
    $ws.Range("CY6").Value = 0.54416761041899997
    $ws.Range("CZ6").Value = 0.455266138165
    $ws.Range("DA6").Value = 0
    $ws.Range("DB6").Value = 1
    $ws.Range("DC6").Value = "Solved Puzzle#4"
    $ws.Range("DD6").Value = 4648
    $ws.Range("DE6").Value = 0.16713271671300001
    $ws.Range("DF6").Value = 1
    $ws.Range("DG6").Value = 3844

# Make "Five Puzzles" the active sheet/tab and set the selection like the
# author left it after entering the new data.
$ws.Activate()
$ws.Range("B5").Select()
